# "NavidGloVe" sheet actually reports Word2Vec results -> fix the sheet name.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NavidGloVe")
$ws.Name = "NavidWord2vec"

# Add the new WEAT 10 result as the new first data row (row 4), pushing the
# existing WEAT 6 / WEAT 7 / WEAT 8 rows down by one.
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "WEAT 10"
$ws.Range("B4").Value = 0.72
$ws.Range("C4").Value = 0.01
$ws.Range("D4").Value = 0.21

# Append the new WEAT 5 result as a new last row.
$ws.Range("A8").Value = "WEAT 5"
$ws.Range("B8").Value = -0.08
$ws.Range("C8").Value = 0.57
$ws.Range("D8").Value = -0.318

# Leave the selection where the author last clicked while editing.
$ws.Range("C11").Select()
